$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 104
$ws.Range("F5").Value = 63
$ws.Range("H5").Value = 63

$ws.Range("E10").Value = 349

$ws.Range("E12").Value = 347
$ws.Range("F12").Value = 192
$ws.Range("H12").Value = 192

$ws.Range("E15").Value = 118

$ws.Range("E22").Value = 135

$ws.Range("E23").Value = 157

$ws.Range("E27").Value = 239

$ws.Range("E28").Value = 145

$ws.Range("E29").Value = 138

$ws.Range("E32").Value = 139
$ws.Range("F32").Value = 75
$ws.Range("H32").Value = 75

$ws.Range("E33").Value = 229

$ws.Range("E34").Value = 162

$ws.Range("E35").Value = 107

$ws.Range("E40").Value = 203

$ws.Range("E41").Value = 293

$ws.Range("E42").Value = 262
$ws.Range("F42").Value = 137
$ws.Range("H42").Value = 137

$ws.Range("F45").Value = 43
$ws.Range("H45").Value = 43

$ws.Range("E46").Value = 228
$ws.Range("F46").Value = 121
$ws.Range("H46").Value = 121

$ws.Range("E47").Value = 329
$ws.Range("F47").Value = 157
$ws.Range("H47").Value = 157

$ws.Range("E48").Value = 149

$ws.Range("E49").Value = 218

$ws.Range("E50").Value = 187
$ws.Range("F50").Value = 66
$ws.Range("H50").Value = 66

$ws.Range("E51").Value = 180
$ws.Range("F51").Value = 69
$ws.Range("H51").Value = 69
